$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "23.394.30"
$ws.Range("E2").Value = "  -0.40%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.625.06"
$ws.Range("E3").Value = "  -0.52%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.000"
$ws.Range("E4").Value = "  -0.13%  "

$ws.Range("E5").Value = "  -0.28%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "303.76"
$ws.Range("E6").Value = "  -1.47%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3786"
$ws.Range("E7").Value = "  +0.38%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "51.99"
$ws.Range("E8").Value = "  -1.51%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3630"
$ws.Range("E9").Value = "  -1.51%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.226"
$ws.Range("E10").Value = "  -3.57%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08095"
$ws.Range("E11").Value = "  -1.08%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.000"
$ws.Range("E12").Value = "  -0.41%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "22.63"
$ws.Range("E13").Value = "  -2.13%  "

$ws.Range("E14").Value = "  -1.63%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.00001246"
$ws.Range("E15").Value = "  -2.62%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "7.212"
$ws.Range("E16").Value = "  -3.25%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.622.70"
$ws.Range("E17").Value = "  -0.98%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "93.48"
$ws.Range("E18").Value = "  -1.38%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06903"
$ws.Range("E19").Value = "  -0.66%  "

$ws.Range("E20").Value = "  -2.72%  "

$ws.Range("E21").Value = "  -0.37%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.406"
$ws.Range("E22").Value = "  -2.57%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "23.397.25"
$ws.Range("E23").Value = "  -0.36%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.71"
$ws.Range("E24").Value = "  -1.93%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.225"
$ws.Range("E25").Value = "  +3.59%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.439"
$ws.Range("E26").Value = "  +1.12%  "

$ws.Range("E27").Value = "  -1.24%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "148.48"
$ws.Range("E28").Value = "  -1.92%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.289"
$ws.Range("E29").Value = "  -0.84%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "134.29"
$ws.Range("E30").Value = "  -1.52%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.303"
$ws.Range("E31").Value = "  -5.01%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.801.89"
$ws.Range("E32").Value = "  -0.69%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.787"
$ws.Range("E33").Value = "  +0.11%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "11.01"
$ws.Range("E34").Value = "  +5.30%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.9521"
$ws.Range("E35").Value = "  -2.14%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02795"
$ws.Range("E36").Value = "  -0.64%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.2528"
$ws.Range("E37").Value = "  -0.30%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.08812"
$ws.Range("E38").Value = "  -0.19%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.104"
$ws.Range("E39").Value = "  -1.77%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.07194"
$ws.Range("E40").Value = "  -3.06%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.356"
$ws.Range("E41").Value = "  -2.63%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.7067"
$ws.Range("E42").Value = "  -1.28%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "16.21"
$ws.Range("E43").Value = "  +0.21%  "

$ws.Range("E44").Value = "  -2.31%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6462"
$ws.Range("E45").Value = "  -2.18%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.325"
$ws.Range("E46").Value = "  -1.20%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.9991"
$ws.Range("E47").Value = "  -0.32%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.989"
$ws.Range("E48").Value = "  -1.45%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.07990"
$ws.Range("E49").Value = "  -0.65%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.203"
$ws.Range("E50").Value = "  -0.99%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "125.51"
$ws.Range("E51").Value = "  -4.23%  "
